$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4
$ws.Range("AF4").Value = 9.5
$ws.Range("G4").Value = 3.3
$ws.Range("I4").Value = 2.25
$ws.Range("T4").Value = 7.5
$ws.Range("W4").Value = 34

# Row 6
$ws.Range("AA6").Value = 6.1
$ws.Range("AB6").Value = 17.5
$ws.Range("AE6").Value = 9.5
$ws.Range("AF6").Value = 22
$ws.Range("AG6").Value = 15
$ws.Range("AH6").Value = 75
$ws.Range("AI6").Value = 50
$ws.Range("AJ6").Value = 65
$ws.Range("G6").Value = 1.88
$ws.Range("H6").Value = 3.05
$ws.Range("I6").Value = 4.3
$ws.Range("L6").Value = 1.42
$ws.Range("M6").Value = 2.47
$ws.Range("N6").Value = 2.2
$ws.Range("O6").Value = 1.52
$ws.Range("P6").Value = 1.47
$ws.Range("Q6").Value = 2.32
$ws.Range("R6").Value = 1.98
$ws.Range("T6").Value = 5.7
$ws.Range("U6").Value = 8
$ws.Range("W6").Value = 16
$ws.Range("X6").Value = 17
$ws.Range("Z6").Value = 7

# Row 7
$ws.Range("AA7").Value = 8.25
$ws.Range("AB7").Value = 25
$ws.Range("AE7").Value = 14.5
$ws.Range("AF7").Value = 40
$ws.Range("AG7").Value = 23
$ws.Range("AH7").Value = 175
$ws.Range("AI7").Value = 90
$ws.Range("AJ7").Value = 100
$ws.Range("G7").Value = 1.44
$ws.Range("H7").Value = 4.05
$ws.Range("I7").Value = 6.7
$ws.Range("M7").Value = 2.82
$ws.Range("R7").Value = 2.18
$ws.Range("S7").Value = 1.53
$ws.Range("U7").Value = 5.8
$ws.Range("W7").Value = 9
$ws.Range("Z7").Value = 9

# Row 8
$ws.Range("AA8").Value = 7.2
$ws.Range("AC8").Value = 120
$ws.Range("AE8").Value = 5.4
$ws.Range("AF8").Value = 6.3
$ws.Range("AI8").Value = 14
$ws.Range("AJ8").Value = 35
$ws.Range("G8").Value = 6.6
$ws.Range("H8").Value = 3.6
$ws.Range("I8").Value = 1.53
$ws.Range("J8").Value = 1.08
$ws.Range("L8").Value = 1.39
$ws.Range("M8").Value = 2.77
$ws.Range("N8").Value = 2.12
$ws.Range("O8").Value = 1.65
$ws.Range("P8").Value = 1.44
$ws.Range("Q8").Value = 2.6
$ws.Range("R8").Value = 2.18
$ws.Range("T8").Value = 13
$ws.Range("U8").Value = 37
$ws.Range("V8").Value = 22
$ws.Range("W8").Value = 175
$ws.Range("Y8").Value = 100

# Row 9
$ws.Range("AE9").Value = 5.1
$ws.Range("AH9").Value = 9.75
$ws.Range("G9").Value = 7.2
$ws.Range("H9").Value = 3.7
$ws.Range("I9").Value = 1.5
$ws.Range("L9").Value = 1.38
$ws.Range("M9").Value = 2.82
$ws.Range("N9").Value = 2.12
$ws.Range("P9").Value = 1.42
$ws.Range("Q9").Value = 2.65
$ws.Range("T9").Value = 15

# Row 10
$ws.Range("AA10").Value = 7.5
$ws.Range("AH10").Value = 15
$ws.Range("G10").Value = 3.75
$ws.Range("H10").Value = 3.8
$ws.Range("I10").Value = 1.83
$ws.Range("R10").Value = 1.75
$ws.Range("S10").Value = 2
$ws.Range("T10").Value = 12
$ws.Range("U10").Value = 21
$ws.Range("V10").Value = 13

# Row 13
$ws.Range("AF13").Value = 51
$ws.Range("N13").Value = 1.44
$ws.Range("O13").Value = 2.7
$ws.Range("P13").Value = 1.19
$ws.Range("Y13").Value = 23

# Row 14
$ws.Range("P14").Value = 1.33

# Row 16
$ws.Range("AA16").Value = 5.9
$ws.Range("AC16").Value = 120
$ws.Range("AE16").Value = 7.1
$ws.Range("AF16").Value = 14.5
$ws.Range("AG16").Value = 11.75
$ws.Range("AH16").Value = 40
$ws.Range("AI16").Value = 32
$ws.Range("AJ16").Value = 50
$ws.Range("G16").Value = 2.62
$ws.Range("H16").Value = 2.87
$ws.Range("I16").Value = 2.92
$ws.Range("J16").Value = 1.12
$ws.Range("K16").Value = 5.7
$ws.Range("L16").Value = 1.5
$ws.Range("M16").Value = 2.47
$ws.Range("N16").Value = 2.47
$ws.Range("O16").Value = 1.5
$ws.Range("P16").Value = 1.55
$ws.Range("Q16").Value = 2.35
$ws.Range("R16").Value = 2
$ws.Range("S16").Value = 1.72
$ws.Range("T16").Value = 6.5
$ws.Range("U16").Value = 12.5
$ws.Range("V16").Value = 11
$ws.Range("W16").Value = 32
$ws.Range("X16").Value = 29
$ws.Range("Z16").Value = 5.7

# Row 18
$ws.Range("AB18").Value = 18
$ws.Range("AE18").Value = 11.75
$ws.Range("AF18").Value = 27
$ws.Range("AH18").Value = 90
$ws.Range("AI18").Value = 55
$ws.Range("G18").Value = 1.65
$ws.Range("H18").Value = 3.55
$ws.Range("I18").Value = 4.85
$ws.Range("M18").Value = 2.77
$ws.Range("N18").Value = 1.98
$ws.Range("O18").Value = 1.65
$ws.Range("P18").Value = 1.44
$ws.Range("Q18").Value = 2.42
$ws.Range("R18").Value = 1.93
$ws.Range("S18").Value = 1.7
$ws.Range("T18").Value = 5.9
$ws.Range("U18").Value = 7.1
$ws.Range("W18").Value = 12.5
$ws.Range("X18").Value = 14.5
$ws.Range("Y18").Value = 32

# Row 23
$ws.Range("AA23").Value = 7.2
$ws.Range("AE23").Value = 9.25
$ws.Range("AF23").Value = 10.5
$ws.Range("AH23").Value = 17.5
$ws.Range("AI23").Value = 13.5
$ws.Range("G23").Value = 3.55
$ws.Range("H23").Value = 3.65
$ws.Range("K23").Value = 8.5
$ws.Range("L23").Value = 1.21
$ws.Range("O23").Value = 2.12
$ws.Range("R23").Value = 1.57
$ws.Range("S23").Value = 2.25
$ws.Range("T23").Value = 13
$ws.Range("U23").Value = 21
$ws.Range("X23").Value = 29
$ws.Range("Z23").Value = 8.5
